$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- D3: the Hydrogen/Non-metallic-minerals cell no longer holds a number;
#     it becomes a blank text cell (matches the other blank cells in the sheet).
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# --- C4: corrected Methanol/Chemicals figure
$ws.Range("C4").Value = 54.40337579526472

# --- C5: corrected Ammonia/Chemicals figure
$ws.Range("C5").Value = 4153.88996024186

# --- Row 7 gets relabelled "Biogas" and now carries a value in column D
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1051.36775377502

# --- New row 8: re-introduce the "Other" row below Biogas, same layout as
#     the former row 7 (bold/bordered label, blank B/C, numeric D = 0)
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 0
